{"js": "// Fix pandoc solution's \"Image Caption\" paragraph style:\n//  - center-align the caption\n//  - use \u9ed1\u4f53 (SimHei) for East-Asian glyphs\n//  - drop the italic styling inherited from the base \"caption\" style\n//  - bump the size to 10.5pt (half-point value 21)\n// and reduce the sample \"Image Caption\" paragraph's first-line indent\n// from 480 twips (24pt) to 420 twips (21pt).\n\nconst styles = context.document.getStyles();\nconst imageCaptionStyle = styles.getByNameOrNullObject(\"Image Caption\");\nawait context.sync();\n\nimageCaptionStyle.paragraphFormat.alignment = Word.Alignment.centered;\nimageCaptionStyle.font.italic = false;\nimageCaptionStyle.font.size = 10.5;\nimageCaptionStyle.font.nameFarEast = \"\u9ed1\u4f53\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.style === \"Image Caption\") {\n    paragraph.paragraphFormat.firstLineIndent = 21; // points (420 twips)\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix pandoc solution's \"Image Caption\" paragraph style:\n#  - center-align the caption\n#  - use \u9ed1\u4f53 (SimHei) for East-Asian glyphs\n#  - drop the italic styling inherited from the base \"caption\" style\n#  - bump the size to 10.5pt (half-point value 21)\n# and reduce the sample \"Image Caption\" paragraph's first-line indent\n# from 480 twips (24pt) to 420 twips (21pt).\n\n$d = $word.ActiveDocument\n\n$imageCaptionStyle = $d.Styles(\"Image Caption\")\n$imageCaptionStyle.ParagraphFormat.Alignment = 1 # wdAlignParagraphCenter\n$imageCaptionStyle.Font.Italic = 0\n$imageCaptionStyle.Font.Size = 10.5\n$imageCaptionStyle.Font.NameFarEast = \"\u9ed1\u4f53\"\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Style.NameLocal -eq \"Image Caption\") {\n        $p.Range.ParagraphFormat.FirstLineIndent = 21 # points (420 twips)\n    }\n}\n"}
